$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customers")
Write-Host $ws.Name
